# Quarterly indexing esoteric bug-fix operation
#
# Each data row (rows 2-16, one per forecast-origin quarter) currently
# stores its q/q error series left-aligned starting in column B. The
# series was mis-indexed by one quarter: every existing value needs to
# shift one column to the right (B->C, C->D, ..., J->K), the value that
# used to sit in the last occupied column falls off the end of the
# row (columns never go past K), and a new value - the q/q error for
# the newly-inserted leading quarter - is written into column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> how many of columns B..K (2..11) are currently populated
# in that row (this shrinks by one every row since each origin quarter
# has one fewer observed error than the quarter before it).
$filledCounts = @{
    2  = 10
    3  = 10
    4  = 10
    5  = 10
    6  = 10
    7  = 9
    8  = 8
    9  = 7
    10 = 6
    11 = 5
    12 = 4
    13 = 3
    14 = 2
    15 = 1
    16 = 0
}

# Row number -> the new value to insert into column B for that row.
$newFirstValues = @{
    2  = -0.6603092772102132
    3  = -0.15162438770796
    4  = -0.2053460154962278
    5  = 0.6162032393936197
    6  = 1.652643173475852
    7  = 0.3110387314724781
    8  = 0.2388379152847414
    9  = 0.6508000635779043
    10 = 0.2387740594105157
    11 = 0.3465902496671606
    12 = 0.00230005330798793
    13 = -0.1902738424076751
    14 = -0.3325070745318338
    15 = 0.1656141382254278
    16 = -0.09587373626955231
}

for ($row = 2; $row -le 16; $row++) {

    $count = $filledCounts[$row]

    # Snapshot the existing B..K values for this row before overwriting
    # anything (column B is index 2, column K is index 11).
    $oldValues = @()
    for ($i = 0; $i -lt $count; $i++) {
        $oldValues += $ws.Cells.Item($row, 2 + $i).Value2
    }

    # Shift everything one column to the right, writing from the
    # rightmost column back towards B so we never clobber a value we
    # still need to read. Anything that would land past column K (11)
    # is simply dropped.
    for ($i = $count - 1; $i -ge 0; $i--) {
        $destCol = 2 + $i + 1
        if ($destCol -le 11) {
            $ws.Cells.Item($row, $destCol).Value2 = $oldValues[$i]
        }
    }

    # Finally, write the new leading value into column B.
    $ws.Cells.Item($row, 2).Value2 = $newFirstValues[$row]
}

Write-Output "quarterly reindex applied to rows 2-16"
